$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/borders/number formats) of column J (rows 2-9) onto the
# new column K so the extra 2021 column visually matches the existing year columns.
$ws.Range("J2:J9").Copy() | Out-Null
$ws.Range("K2:K9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new 2021 figures in column K.
$ws.Range("K3").Value = 2021
$ws.Range("K4").Value = 295
$ws.Range("K5").Value = 163
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 27
$ws.Range("K9").Value = 8

# Move the active selection, matching the saved sheet view.
$ws.Range("L5").Select() | Out-Null
